$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14 through 24 (inclusive) - these are the "middle" province
# rows that got removed; row 25 (id_kho 22168000, Hoai Nhon-Binh Dinh)
# shifts up to become the new row 14.
$ws.Rows("14:24").Delete()

# The ten_kho/ghi_chu columns (B, F) and the trailing data row's id_kho
# cell (A14) no longer carry an explicit cell style in the edited file.
$ws.Range("B2:B14").Style = "Normal"
$ws.Range("F2:F14").Style = "Normal"
$ws.Range("A14").Style = "Normal"

# Update the active selection to match the edited workbook.
$ws.Range("E21").Select()
